# Remove the trailing "Ver no Jupiter..." / "© 2020 ..." footer block
# (and the blank paragraph immediately preceding it) that used to follow
# the "LOQ4233: Gestão de Negócios (Requisito fraco)" requirement line.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter*") {
        # the blank paragraph right before this one is also removed
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -like "*Powered by Jekyll*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
